$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely (it duplicated column F's GENE values and carried the bold/bordered style);
# this shifts columns B:F left to A:E, so the "MODEL_CONDITION" header ends up in column D.
$ws.Columns("A").Delete()

# Rename the "MODEL_CONDITION" header text to "MODELCONDITION".
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")
